$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 25000
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H36").Value = 25000
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H112").Value = 2632.2942
$ws.Range("I112").Value = 5000
$ws.Range("J112").Value = 2484.3125
$ws.Range("K112").Value = 15000
$ws.Range("L112").Value = 7452.9375
$ws.Range("M112").Value = -13892
$ws.Range("N112").Value = -9668.9375

$ws.Range("H129").Value = 1753.5834
$ws.Range("J129").Value = 1994.6666
$ws.Range("L129").Value = 5983.9998
$ws.Range("N129").Value = -15983.9998

$ws.Range("H132").Value = 1165.9231
$ws.Range("I132").Value = 1165.9231
$ws.Range("K132").Value = 3497.7693
$ws.Range("M132").Value = -967.7692999999999

$ws.Range("H135").Value = 7172.76
$ws.Range("J135").Value = 10729.857
$ws.Range("L135").Value = 96568.713
$ws.Range("N135").Value = -101638.713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1709.625
$ws.Range("I45").Value = 730.36365
$ws.Range("K45").Value = 730.36365
$ws.Range("M45").Value = -353.36365

$ws.Range("H115").Value = 80593.336
$ws.Range("J115").Value = 80593.336
$ws.Range("L115").Value = 80593.336
$ws.Range("N115").Value = -83727.336

$ws.Range("H132").Value = 6115.4136
$ws.Range("I132").Value = 3464.4348
$ws.Range("J132").Value = 16277.5
$ws.Range("K132").Value = 10393.3044
$ws.Range("L132").Value = 48832.5
$ws.Range("M132").Value = -7863.304400000001
$ws.Range("N132").Value = -53892.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H127").Value = 60474.25
$ws.Range("J127").Value = 60474.25
$ws.Range("L127").Value = 60474.25
$ws.Range("N127").Value = -70394.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 895.3333
$ws.Range("I22").Value = 868.5
$ws.Range("K22").Value = 868.5
$ws.Range("M22").Value = -518.5

$ws.Range("H28").Value = 100772.6
$ws.Range("J28").Value = 100772.6
$ws.Range("L28").Value = 100772.6
$ws.Range("N28").Value = -101262.6

$ws.Range("H31").Value = 972743.9
$ws.Range("I31").Value = 1388.3529
$ws.Range("K31").Value = 1388.3529
$ws.Range("M31").Value = -1093.3529

$ws.Range("H34").Value = 972743.9
$ws.Range("I34").Value = 1388.3529
$ws.Range("K34").Value = 1388.3529
$ws.Range("M34").Value = -1186.3529

$ws.Range("H117").Value = 39949.5
$ws.Range("J117").Value = 39949.5
$ws.Range("L117").Value = 39949.5
$ws.Range("N117").Value = -49127.5

$ws.Range("H122").Value = 4857.65
$ws.Range("I122").Value = 1973.7142
$ws.Range("K122").Value = 5921.142599999999
$ws.Range("M122").Value = -3471.142599999999

$ws.Range("H129").Value = 60125.4
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 60125.4
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 60125.4
$ws.Range("N129").Value = -70125.39999999999
$ws.Range("M129").ClearContents()

$ws.Range("H134").Value = 558796.5600000001
$ws.Range("I134").Value = 770679.4
$ws.Range("K134").Value = 2312038.2
$ws.Range("M134").Value = -2309503.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1115.8182
$ws.Range("I113").Value = 414.66666
$ws.Range("J113").Value = 1378.75
$ws.Range("K113").Value = 1243.99998
$ws.Range("L113").Value = 4136.25
$ws.Range("M113").Value = 926.0000199999999
$ws.Range("N113").Value = -8476.25

$ws.Range("H134").Value = 4434.933
$ws.Range("I134").Value = 3543.6667
$ws.Range("K134").Value = 10631.0001
$ws.Range("M134").Value = -5561.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 75000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H80").Value = 9356.261
$ws.Range("I80").Value = 5344.222
$ws.Range("J80").Value = 11935.429
$ws.Range("K80").Value = 5344.222
$ws.Range("L80").Value = 11935.429
$ws.Range("M80").Value = -4346.222
$ws.Range("N80").Value = -13931.429

$ws.Range("H83").Value = 9356.261
$ws.Range("I83").Value = 5344.222
$ws.Range("J83").Value = 11935.429
$ws.Range("K83").Value = 26721.11
$ws.Range("L83").Value = 59677.145
$ws.Range("M83").Value = -21729.11
$ws.Range("N83").Value = -69661.145

$ws.Range("H104").Value = 66412.75
$ws.Range("J104").Value = 66412.75
$ws.Range("L104").Value = 66412.75
$ws.Range("N104").Value = -73400.75

$ws.Range("H111").Value = 47066.332
$ws.Range("J111").Value = 47066.332
$ws.Range("L111").Value = 47066.332
$ws.Range("N111").Value = -53200.332

$ws.Range("H116").Value = 69995
$ws.Range("J116").Value = 69995
$ws.Range("L116").Value = 69995
$ws.Range("N116").Value = -79173

$ws.Range("H118").Value = 43619.8
$ws.Range("J118").Value = 43619.8
$ws.Range("L118").Value = 43619.8
$ws.Range("N118").Value = -46933.8

$ws.Range("H122").Value = 3848.6
$ws.Range("I122").Value = 3060.875
$ws.Range("K122").Value = 9182.625
$ws.Range("M122").Value = -6732.625

$ws.Range("H128").Value = 110970
$ws.Range("J128").Value = 110970
$ws.Range("L128").Value = 110970
$ws.Range("N128").Value = -120930

$ws.Range("H129").Value = 69161.664
$ws.Range("J129").Value = 69161.664
$ws.Range("L129").Value = 69161.664
$ws.Range("N129").Value = -79161.664

$ws.Range("H130").Value = 85997
$ws.Range("J130").Value = 85997
$ws.Range("L130").Value = 85997
$ws.Range("N130").Value = -96037

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 50000000
$ws.Range("J4").Value = 50000000
$ws.Range("L4").Value = 50000000
$ws.Range("N4").Value = -50000226

$ws.Range("H28").Value = 50000000
$ws.Range("J28").Value = 50000000
$ws.Range("L28").Value = 50000000
$ws.Range("N28").Value = -50000464

$ws.Range("H37").Value = 50000000
$ws.Range("J37").Value = 50000000
$ws.Range("L37").Value = 50000000
$ws.Range("N37").Value = -50000214

$ws.Range("H47").Value = 17831.666
$ws.Range("J47").Value = 34495
$ws.Range("L47").Value = 34495
$ws.Range("N47").Value = -35475

$ws.Range("H52").Value = 17831.666
$ws.Range("J52").Value = 34495
$ws.Range("L52").Value = 34495
$ws.Range("N52").Value = -34961

$ws.Range("H74").Value = 116000
$ws.Range("J74").Value = 116000
$ws.Range("L74").Value = 116000
$ws.Range("N74").Value = -117996

$ws.Range("H77").Value = 116000
$ws.Range("J77").Value = 116000
$ws.Range("L77").Value = 348000
$ws.Range("N77").Value = -357984

$ws.Range("H100").Value = 5030.625
$ws.Range("I100").Value = 5874.1665
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 5874.1665
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -5333.1665
$ws.Range("N100").Value = -3582

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 75000
$ws.Range("J21").Value = 75000
$ws.Range("L21").Value = 75000
$ws.Range("N21").Value = -75470

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H35").Value = 75000
$ws.Range("J35").Value = 75000
$ws.Range("L35").Value = 75000
$ws.Range("N35").Value = -75580

$ws.Range("H81").Value = 34464.668
$ws.Range("I81").Value = 1295
$ws.Range("K81").Value = 2590
$ws.Range("M81").Value = -1529

$ws.Range("H84").Value = 34464.668
$ws.Range("I84").Value = 1295
$ws.Range("K84").Value = 12950
$ws.Range("M84").Value = -7646

$ws.Range("H141").Value = 61996
$ws.Range("J141").Value = 61996
$ws.Range("L141").Value = 61996
$ws.Range("N141").Value = -72356
